$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.985.77'
$ws.Range('E2').Value = '  +0.24%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.638.52'
$ws.Range('E3').Value = '  -0.62%  '

$ws.Range('E4').Value = '  -0.51%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.14'
$ws.Range('E5').Value = '  -0.33%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5131'
$ws.Range('E6').Value = '  +0.37%  '

$ws.Range('E7').Value = '  -0.54%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2577'
$ws.Range('E8').Value = '  -0.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06348'
$ws.Range('E9').Value = '  -1.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.80'
$ws.Range('E10').Value = '  +0.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07765'
$ws.Range('E11').Value = '  -0.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.274'
$ws.Range('E12').Value = '  -1.10%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.635.70'
$ws.Range('E13').Value = '  -2.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5462'
$ws.Range('E14').Value = '  -0.33%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅7757'
$ws.Range('E15').Value = '  -1.89%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.35'
$ws.Range('E16').Value = '  -1.09%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.988.43'
$ws.Range('E17').Value = '  -0.11%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9972'
$ws.Range('E18').Value = '  -0.98%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '197.54'
$ws.Range('E19').Value = '  -0.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.443'
$ws.Range('E20').Value = '  +0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.929'
$ws.Range('E21').Value = '  -1.18%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.085'
$ws.Range('E22').Value = '  +0.10%  '

$ws.Range('E23').Value = '  -0.74%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.928'
$ws.Range('E24').Value = '  +3.50%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.88'
$ws.Range('E25').Value = '  +0.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1232'
$ws.Range('E26').Value = '  +7.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.845'
$ws.Range('E27').Value = '  -0.93%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.58'
$ws.Range('E28').Value = '  -1.30%  '

$ws.Range('E29').Value = '  -0.06%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04849'
$ws.Range('E30').Value = '  -3.34%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.287'
$ws.Range('E31').Value = '  +0.20%  '

$ws.Range('E32').Value = '  +0.26%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.539'
$ws.Range('E33').Value = '  -0.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9141'
$ws.Range('E35').Value = '  +1.95%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.567'
$ws.Range('E36').Value = '  -0.97%  '

$ws.Range('E37').Value = '  -0.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.124.64'
$ws.Range('E38').Value = '  -0.75%  '

$ws.Range('E39').Value = '  +0.19%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').Value = '  -0.73%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.522'
$ws.Range('E41').Value = '  -1.83%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.568'
$ws.Range('E42').Value = '  -1.87%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8072'
$ws.Range('E43').Value = '  -1.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.34'
$ws.Range('E44').Value = '  -0.49%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₈121'
$ws.Range('E45').Value = '  -2.62%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.780.33'
$ws.Range('E46').Value = '  -0.29%  '

$ws.Range('E47').Value = '  -0.17%  '

$ws.Range('E48').Value = '  -0.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  -0.59%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05215'
$ws.Range('E50').Value = '  +2.34%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.473'
$ws.Range('E51').Value = '  +0.42%  '
